$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet "About" (1st sheet): rebuild content for the new
# large/small primary-energy-unit layout.
# ============================================================
$ws1 = $wb.Worksheets.Item(1)

# Drop the old hyperlink (and its "Hyperlink" cell style) since
# the new About sheet no longer references an external source.
$ws1.Hyperlinks.Delete()
$ws1.Cells.Clear()
$wb.Styles.Item("Hyperlink").Delete()

$ws1.Range("A1").Value = "BpTPEU BTU per Large Primary Energy Unit"
$ws1.Range("A1").Font.Bold = $true

$ws1.Range("A2").Value = "BpTPEU BTU per Small Primary Energy Unit"
$ws1.Range("A2").Font.Bold = $true

$ws1.Range("A4").Value = "Source:"
$ws1.Range("A4").Font.Bold = $true
$ws1.Range("B4").Value = "none needed"

$ws1.Range("B5").HorizontalAlignment = -4131   # xlLeft

$ws1.Range("A9").Value = "Notes"
$ws1.Range("A9").Font.Bold = $true

$ws1.Range("A10").Value = "For the U.S.:"
$ws1.Range("A11").Value = "The large primary energy output unit (used in totals graphs) is: quadrillion BTU"
$ws1.Range("A12").Value = "The small primary energy output unit (used in energy intensity per unit GDP graphs) is: thousand BTU"

$ws1.PageSetup.Orientation = 1   # xlPortrait

# ============================================================
# Sheet "BpTPEU" -> "BpTPEU-large": the large-unit conversion
# sheet (quadrillion BTU, i.e. 10^15).
# ============================================================
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "BpTPEU-large"

$ws2.Range("B1").Value = "large primary energy output unit"
$ws2.Range("B1").HorizontalAlignment = -4131   # xlLeft

$ws2.Range("B2").Formula = "=10^15"

# ============================================================
# New sheet "BpTPEU-small": duplicate the large-unit sheet so
# that tab color / column widths / number formats all carry
# over, then edit it to describe the small-unit (thousand BTU,
# i.e. 10^3) conversion.
# ============================================================
$ws2.Copy([Type]::Missing, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "BpTPEU-small"

$ws3.Range("B1").Value = "small primary energy output unit"

$ws3.Range("B2").Formula = "=10^3"
$ws3.Range("B2").NumberFormat = "General"

# Restore the original active sheet/tab selection.
$ws1.Activate()
